$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-09 Saturday" "2024-03-10 Sunday"
Replace-Text "240÷5=" "705÷7="
Replace-Text "532÷8=" "814÷5="
Replace-Text "826÷3=" "978÷5="
Replace-Text "722÷4=" "834÷4="
Replace-Text "546÷3=" "851÷8="
Replace-Text "522÷7=" "126÷2="
Replace-Text "497÷7=" "906÷4="
Replace-Text "156÷5=" "918÷7="
Replace-Text "659÷5=" "794÷4="
Replace-Text "191÷4=" "630÷2="
Replace-Text "196÷9=" "428÷8="
Replace-Text "363÷2=" "173÷8="
Replace-Text "357÷9=" "394÷4="
Replace-Text "377÷7=" "135÷6="
Replace-Text "962÷5=" "629÷9="
Replace-Text "709÷9=" "101÷4="
Replace-Text "644÷8=" "411÷3="
Replace-Text "312÷8=" "843÷2="
Replace-Text "525÷3=" "350÷2="
Replace-Text "770÷3=" "845÷4="
Replace-Text "766÷8=" "376÷7="
Replace-Text "374÷7=" "880÷5="
Replace-Text "592÷7=" "840÷2="
Replace-Text "492÷9=" "658÷5="
Replace-Text "321÷5=" "826÷6="
